$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "97.700.90"
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  +0.87%  "
$cell.Style = "Normal"

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.298.41"
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  -0.84%  "
$cell.Style = "Normal"

# Row 4
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  -0.10%  "
$cell.Style = "Normal"

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "254.93"
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  +4.56%  "
$cell.Style = "Normal"

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "619.49"
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  -0.43%  "
$cell.Style = "Normal"

# Row 7
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "1.42"
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  +26.14%  "
$cell.Style = "Normal"

# Row 8
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.398"
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  +2.94%  "
$cell.Style = "Normal"

# Row 9
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  -0.01%  "
$cell.Style = "Normal"

# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.892"
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  +13.50%  "
$cell.Style = "Normal"

# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "3.301.06"
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  -0.86%  "
$cell.Style = "Normal"

# Row 12
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  -0.65%  "
$cell.Style = "Normal"

# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "38.22"
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  +8.58%  "
$cell.Style = "Normal"

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "97.358.45"
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  +0.53%  "
$cell.Style = "Normal"

# Row 15
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  +0.34%  "
$cell.Style = "Normal"

# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "3.917.01"
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  -1.17%  "
$cell.Style = "Normal"

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "5.44"
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  -0.70%  "
$cell.Style = "Normal"

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "3.294.01"
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  -1.91%  "
$cell.Style = "Normal"

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "3.51"
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  -1.65%  "
$cell.Style = "Normal"

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "15.05"
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  -0.49%  "
$cell.Style = "Normal"

# Row 21
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "6.09"
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  +4.33%  "
$cell.Style = "Normal"

# Row 22
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "477.20"
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  -1.96%  "
$cell.Style = "Normal"

# Row 23
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "9.37"
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  +1.29%  "
$cell.Style = "Normal"

# Row 24
$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "0.0000203"
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  -2.20%  "
$cell.Style = "Normal"

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "5.57"
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  -1.63%  "
$cell.Style = "Normal"

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "87.58"
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  -1.08%  "
$cell.Style = "Normal"

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "11.78"
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "  -2.55%  "
$cell.Style = "Normal"

# Row 28
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "3.471.73"
$cell.Style = "Normal"
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  -1.86%  "
$cell.Style = "Normal"

# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "0.290"
$cell.Style = "Normal"
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  +20.90%  "
$cell.Style = "Normal"

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = "Normal"
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  -0.28%  "
$cell.Style = "Normal"

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.186"
$cell.Style = "Normal"
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  +2.89%  "
$cell.Style = "Normal"

# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.130"
$cell.Style = "Normal"
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  +7.65%  "
$cell.Style = "Normal"

# Row 33
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "9.72"
$cell.Style = "Normal"
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  +4.94%  "
$cell.Style = "Normal"

# Row 34
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.997"
$cell.Style = "Normal"
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  -0.37%  "
$cell.Style = "Normal"

# Row 35
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "27.46"
$cell.Style = "Normal"
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  +0.07%  "
$cell.Style = "Normal"

# Row 36
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "7.13"
$cell.Style = "Normal"
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  -3.19%  "
$cell.Style = "Normal"

# Row 37
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.146"
$cell.Style = "Normal"
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  -2.95%  "
$cell.Style = "Normal"

# Row 38
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  -0.51%  "
$cell.Style = "Normal"

# Row 39
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "24.79"
$cell.Style = "Normal"
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  +0.70%  "
$cell.Style = "Normal"

# Row 40
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.453"
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  +1.11%  "
$cell.Style = "Normal"

# Row 41
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "487.82"
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  -0.93%  "
$cell.Style = "Normal"

# Row 42
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  +6.22%  "
$cell.Style = "Normal"

# Row 43
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  -2.39%  "
$cell.Style = "Normal"

# Row 44
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.796"
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  -0.44%  "
$cell.Style = "Normal"

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "3.14"
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  -2.98%  "
$cell.Style = "Normal"

# Row 47
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "158.53"
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  -0.90%  "
$cell.Style = "Normal"

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.89"
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  -1.52%  "
$cell.Style = "Normal"

# Row 49
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.837"
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  +4.25%  "
$cell.Style = "Normal"

# Row 50
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "4.58"
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  +1.03%  "
$cell.Style = "Normal"

# Row 51
$ws.Range("B51").Value = "Optimism"
$ws.Range("C51").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "2.17"
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  +3.05%  "
$cell.Style = "Normal"
